# Added Framework to GetTableCellTextTest
# Appends extra test rows (including some autocomplete-duplicated values,
# matching the original author's manual data entry) to the Company,
# Contact and Country lookup sheets used by the table-cell-text test.

$wb = $excel.ActiveWorkbook

# --- Company sheet: extend A8:A14 -------------------------------------
$wsCompany = $wb.Worksheets.Item("Company")
$wsCompany.Range("A8").Value2  = "Eliran Duveen"
$wsCompany.Range("A9").Value2  = "sdfdh"
$wsCompany.Range("A10").Value2 = "h5f23"
$wsCompany.Range("A11").Value2 = "sdgvr dsrgr"
$wsCompany.Range("A12").Value2 = "Ernst Handel"
$wsCompany.Range("A12").NumberFormat = "General"
$wsCompany.Range("A13").Value2 = "/.,/rty5"
$wsCompany.Range("A14").Value2 = "Laughing Bacchus Winecellars"
$wsCompany.Range("A14").NumberFormat = "General"
$wsCompany.Range("A14").Select() | Out-Null

# --- Contact sheet: extend A8:A14 --------------------------------------
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Range("A8").Value2  = "Eliran Duveen"
$wsContact.Range("A9").Value2  = "sgdsgh"
$wsContact.Range("A10").Value2 = "''lih[oi"
$wsContact.Range("A11").Value2 = "cs;kljdb"
$wsContact.Range("A12").Value2 = "f234thh6;"
$wsContact.Range("A13").Value2 = "Maria Anders"
$wsContact.Range("A14").Value2 = "Yoshi Tannamuri"
$wsContact.Range("A2").Select() | Out-Null

# --- Country sheet: extend A8:A13 --------------------------------------
$wsCountry = $wb.Worksheets.Item("Country")
$wsCountry.Range("A8").Value2  = "Israel"
$wsCountry.Range("A9").Value2  = "Maria Anders"
$wsCountry.Range("A10").Value2 = "gdeg"
$wsCountry.Range("A11").Value2 = "nhfjh345"
$wsCountry.Range("A12").Value2 = "Mexico"
$wsCountry.Range("A13").Value2 = "UK"
$wsCountry.Range("A13").Select() | Out-Null
